$d = $word.ActiveDocument

# --- Recolor the first "Pass application verifier." / "Leak testing." pair
#     (immediately following "Unit tests.") from red (FF0000) to purple (7030A0).
$d.Paragraphs.Item(9).Range.Font.Color = 10498160   # 0x7030A0 -> Pass application verifier.
$d.Paragraphs.Item(10).Range.Font.Color = 10498160  # 0x7030A0 -> Leak testing.

# --- Recolor the second "Pass application verifier." / "Leak testing." pair
#     (further down the Todo list) from red (FF0000) to purple (7030A0).
$d.Paragraphs.Item(17).Range.Font.Color = 10498160  # 0x7030A0 -> Pass application verifier.
$d.Paragraphs.Item(18).Range.Font.Color = 10498160  # 0x7030A0 -> Leak testing.

# --- Move the "_GoBack" bookmark from the end of "Unit tests." to the middle of
#     "... meet the minimum requirements." (splitting the word "requirements.").
#     Adding a bookmark with a name that already exists moves/replaces it, so the
#     old location is cleared automatically.
$lastPara = $d.Paragraphs.Item(25)
$lastText = $lastPara.Range.Text
$splitOffset = $lastText.IndexOf("requirements.") + ("requirem").Length
$splitPos = $lastPara.Range.Start + $splitOffset
$bmRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
